$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 20
$ws1.Range("F7").Value = 291
$ws1.Range("F10").Value = 18
$ws1.Range("F12").Value = 558
$ws1.Range("F14").Value = 13198
$ws1.Range("F17").Value = 15
$ws1.Range("F18").Value = 5429
$ws1.Range("F19").Value = 5561
$ws1.Range("F20").Value = 27

# Sheet "全部类型" (sheet index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 20
$ws4.Range("F23").Value = 291
$ws4.Range("F32").Value = 18
$ws4.Range("F34").Value = 558
$ws4.Range("F36").Value = 13198
$ws4.Range("F40").Value = 15
$ws4.Range("F41").Value = 5429
$ws4.Range("F42").Value = 5561
$ws4.Range("F43").Value = 27
